# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# just before column N ("Late"), pushing the existing "Late",
# "heading" (Paid Date header above Outstanding) and "Outstanding"
# columns one place to the right (N->O, O->P, P->Q). The new column N
# keeps the same width as column M (which stores as width 11 in the
# sheet XML) but carries no header text / no bestFit flag, matching a
# plain "Insert Column" operation in Excel.
#
# The active sheet also changes from "Edit Repayment Schedule" to
# "Repayment schedule", and the selected cell/active cell on each of
# those two sheets changes as well.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsEdit  = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new blank column before column N, carrying column M's width.
$mw = $wsRepay.Columns("M:M").ColumnWidth
$wsRepay.Columns("N:N").Insert() | Out-Null
$wsRepay.Columns("N:N").ColumnWidth = $mw

# Update the stored selection on "Edit Repayment Schedule" (it is no
# longer the active tab, but its own last-selected cell still moves).
$wsEdit.Activate() | Out-Null
$wsEdit.Range("B5").Select() | Out-Null

# Make "Repayment schedule" the active sheet/tab with its new selection.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("I22").Select() | Out-Null
